$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C5) from 45174 to 45175 (one day later)
$ws.Range("C2:C5").Value = 45175
